# UGCOfWorkshopMonth.xlsx update
# - Rolls the "TUS" snapshot date forward (2024-04-09 -> 2024-04-22)
# - Refreshes the Chinese-language leaderboard (rows 3-6) with new entries
# - Refreshes the "Model" leaderboard (previously rows 24-26, now rows 8-12)
#   with new entries, adding a 4th row
# - Shrinks the sheet's used range from A1:I26 down to A1:I12

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($row, $col, $text) {
    # Some values (e.g. "02/03/2024", "0") look like dates/numbers to
    # Excel's input parser and would otherwise get silently coerced into a
    # date serial or a numeric value. Briefly flipping the cell to Text
    # format forces the literal string to be stored, then flipping the
    # style back to Normal drops the formatting override again so the
    # cell ends up with no explicit style (matching plain text cells
    # elsewhere on the sheet).
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

# Start clean: wipe the whole sheet and rebuild it top-to-bottom so the
# rebuilt shared-strings table comes out in natural reading order (this
# also drops every now-unused shared string from the old revision).
$ws.Cells.Clear()

# --- Row 1 ---
$ws.Cells.Item(1, 1).Value = "Level"

# --- Row 2: header row ---
$ws.Cells.Item(2, 1).Value = "Rank"
$ws.Cells.Item(2, 2).Value = "Title"
$ws.Cells.Item(2, 3).Value = "Creator"
$ws.Cells.Item(2, 4).Value = "Date Posted"
$ws.Cells.Item(2, 5).Value = "Country"
$ws.Cells.Item(2, 6).Value = "Language"
$ws.Cells.Item(2, 7).Value = "TUS (2024-04-22)"
$ws.Cells.Item(2, 8).Value = "Rating"
$ws.Cells.Item(2, 9).Value = "Comment Count"

# --- Rows 3-6: Chinese leaderboard entries ---
$ws.Cells.Item(3, 1).Value = 1
$ws.Cells.Item(3, 2).Value = "诀别书"
$ws.Cells.Item(3, 3).Value = "可乐没有气了"
Set-TextValue 3 4 "19/03/2024"
$ws.Cells.Item(3, 5).Value = "N/A"
$ws.Cells.Item(3, 6).Value = "zh-cn"
$ws.Cells.Item(3, 7).Value = 184
$ws.Cells.Item(3, 8).Value = "N/A"
Set-TextValue 3 9 "0"

$ws.Cells.Item(4, 1).Value = 2
$ws.Cells.Item(4, 2).Value = "BOB BALL3"
$ws.Cells.Item(4, 3).Value = "脸红"
Set-TextValue 4 4 "21/03/2024"
$ws.Cells.Item(4, 5).Value = "CN"
$ws.Cells.Item(4, 6).Value = "zh-cn"
$ws.Cells.Item(4, 7).Value = 1518
$ws.Cells.Item(4, 8).Value = "N/A"
Set-TextValue 4 9 "0"

$ws.Cells.Item(5, 1).Value = 3
$ws.Cells.Item(5, 2).Value = "禁止摆烂萌新图"
$ws.Cells.Item(5, 3).Value = "肉女士"
Set-TextValue 5 4 "23/03/2024"
$ws.Cells.Item(5, 5).Value = "N/A"
$ws.Cells.Item(5, 6).Value = "zh-cn"
$ws.Cells.Item(5, 7).Value = 1195
$ws.Cells.Item(5, 8).Value = "N/A"
Set-TextValue 5 9 "0"

$ws.Cells.Item(6, 1).Value = 4
$ws.Cells.Item(6, 2).Value = "PartyTest_Fyang"
$ws.Cells.Item(6, 3).Value = "肥羊"
Set-TextValue 6 4 "27/03/2024"
$ws.Cells.Item(6, 5).Value = "CN"
$ws.Cells.Item(6, 6).Value = "zh-cn"
$ws.Cells.Item(6, 7).Value = 500
$ws.Cells.Item(6, 8).Value = "N/A"
Set-TextValue 6 9 "0"

# --- Row 7: "Model" section label ---
$ws.Cells.Item(7, 1).Value = "Model"

# --- Row 8: header row (repeated) ---
$ws.Cells.Item(8, 1).Value = "Rank"
$ws.Cells.Item(8, 2).Value = "Title"
$ws.Cells.Item(8, 3).Value = "Creator"
$ws.Cells.Item(8, 4).Value = "Date Posted"
$ws.Cells.Item(8, 5).Value = "Country"
$ws.Cells.Item(8, 6).Value = "Language"
$ws.Cells.Item(8, 7).Value = "TUS (2024-04-22)"
$ws.Cells.Item(8, 8).Value = "Rating"
$ws.Cells.Item(8, 9).Value = "Comment Count"

# --- Rows 9-12: "Model" leaderboard entries ---
$ws.Cells.Item(9, 1).Value = 1
$ws.Cells.Item(9, 2).Value = "Dogday"
$ws.Cells.Item(9, 3).Value = "Kimmel"
Set-TextValue 9 4 "23/03/2024"
$ws.Cells.Item(9, 5).Value = "N/A"
$ws.Cells.Item(9, 6).Value = "es"
$ws.Cells.Item(9, 7).Value = 35
$ws.Cells.Item(9, 8).Value = "N/A"
Set-TextValue 9 9 "0"

$ws.Cells.Item(10, 1).Value = 2
$ws.Cells.Item(10, 2).Value = "Adolf Hitler"
$ws.Cells.Item(10, 3).Value = "ebrunedre"
Set-TextValue 10 4 "23/03/2024"
$ws.Cells.Item(10, 5).Value = "TR"
$ws.Cells.Item(10, 6).Value = "tr"
$ws.Cells.Item(10, 7).Value = 456
$ws.Cells.Item(10, 8).Value = "N/A"
Set-TextValue 10 9 "2"

$ws.Cells.Item(11, 1).Value = 3
$ws.Cells.Item(11, 2).Value = "RED ROBIN!"
$ws.Cells.Item(11, 3).Value = "ordinalst"
Set-TextValue 11 4 "07/03/2024"
$ws.Cells.Item(11, 5).Value = "N/A"
$ws.Cells.Item(11, 6).Value = "en"
$ws.Cells.Item(11, 7).Value = 117
$ws.Cells.Item(11, 8).Value = "N/A"
Set-TextValue 11 9 "0"

$ws.Cells.Item(12, 1).Value = 4
$ws.Cells.Item(12, 2).Value = "КЛОУН"
$ws.Cells.Item(12, 3).Value = "7700n"
Set-TextValue 12 4 "02/03/2024"
$ws.Cells.Item(12, 5).Value = "FR"
$ws.Cells.Item(12, 6).Value = "ru"
$ws.Cells.Item(12, 7).Value = 1053
$ws.Cells.Item(12, 8).Value = "N/A"
Set-TextValue 12 9 "0"
